# feat(PYME-4265): Add subscription id in excel
#
# Insert a new "Subscription ID" column before the existing "Antivirus
# Quantity" column (i.e. as the new column I), shifting every column
# from I..M one place to the right (J..N), and keep the header row,
# the table's autofilter and its related defined name in sync.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at position 9 (I); existing I:M shift to J:N.
$ws.Columns.Item(9).Insert()

# Fill in the header for the newly inserted column.
$ws.Cells.Item(1, 9).Value = "Subscription ID"

# Match the width Excel would use for a freshly inserted column here,
# copying it from the neighboring "Customer Tax ID" column (H).
$ws.Columns.Item(9).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# The autofilter range needs to grow from A1:M1 to A1:N1 to include the
# new column. Toggling it off/on forces the stored range to refresh.
$ws.AutoFilterMode = $false
$ws.Range("A1:N1").AutoFilter(1) | Out-Null

# Keep the workbook-level hidden "_FilterDatabase" name (used by the
# autofilter) pointing at the same, now-widened, range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$N`$1"
    }
}

# Restore the active cell selection seen in the edited workbook.
$ws.Range("I2").Select() | Out-Null
